$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 9.596329634041382
$ws.Range("D2").Value = 7.847387925946922
$ws.Range("E2").Value = 11.9118966710953
$ws.Range("F2").Value = 37.87901270813921
$ws.Range("G2").Value = 46.94130587132261
$ws.Range("H2").Value = 17.89236940006412
$ws.Range("I2").Value = 22.711319792529
$ws.Range("J2").Value = 9.307506806713441
$ws.Range("M2").Value = 29.9279915337514
$ws.Range("B3").Value = 9.173258789433474
$ws.Range("D3").Value = 7.810326985520167
$ws.Range("E3").Value = 11.93703004693193
$ws.Range("F3").Value = 37.73774504488436
$ws.Range("G3").Value = 46.25926200953842
$ws.Range("H3").Value = 17.87754352438132
$ws.Range("I3").Value = 22.17744026971472
$ws.Range("J3").Value = 9.358804115323196
$ws.Range("M3").Value = 28.76838881605268
$ws.Range("B4").Value = 8.902469131019926
$ws.Range("D4").Value = 7.78741155923098
$ws.Range("E4").Value = 11.95434399973242
$ws.Range("F4").Value = 37.66856512642453
$ws.Range("G4").Value = 45.86136672078449
$ws.Range("H4").Value = 17.87516887825283
$ws.Range("I4").Value = 21.84971307132503
$ws.Range("J4").Value = 9.391955545518341
$ws.Range("M4").Value = 28.0311666719889
$ws.Range("B5").Value = 8.78946025951371
$ws.Range("D5").Value = 7.778035479297832
$ws.Range("E5").Value = 11.9618736521008
$ws.Range("F5").Value = 37.64478074656034
$ws.Range("G5").Value = 45.70467724289239
$ws.Range("H5").Value = 17.87588356724648
$ws.Range("I5").Value = 21.71640672218503
$ws.Range("J5").Value = 9.405882251024639
$ws.Range("M5").Value = 27.72476649264869
$ws.Range("B6").Value = 8.77053817205319
$ws.Range("D6").Value = 7.776476369590624
$ws.Range("E6").Value = 11.96315260211498
$ws.Range("F6").Value = 37.64109713832853
$ws.Range("G6").Value = 45.67899425202508
$ws.Range("H6").Value = 17.87610353904637
$ws.Range("I6").Value = 21.69429285050338
$ws.Range("J6").Value = 9.408220003314021
$ws.Range("M6").Value = 27.67353959582947
$ws.Range("B7").Value = 8.900955663553747
$ws.Range("D7").Value = 7.787285260359809
$ws.Range("E7").Value = 11.95444362663496
$ws.Range("F7").Value = 37.66822653505635
$ws.Range("G7").Value = 45.85923120218094
$ws.Range("H7").Value = 17.87517171877072
$ws.Range("I7").Value = 21.84791395894061
$ws.Range("J7").Value = 9.392141675048563
$ws.Range("M7").Value = 28.02705813327772
$ws.Range("B8").Value = 9.45280450203725
$ws.Range("D8").Value = 7.834643900507891
$ws.Range("E8").Value = 11.9201725733219
$ws.Range("F8").Value = 37.82665666855829
$ws.Range("G8").Value = 46.70193124954226
$ws.Range("H8").Value = 17.88585570704038
$ws.Range("I8").Value = 22.52735521546734
$ws.Range("J8").Value = 9.324851380234985
$ws.Range("M8").Value = 29.53360128079446
$ws.Range("B9").Value = 10.44341174515645
$ws.Range("D9").Value = 7.926143913846621
$ws.Range("E9").Value = 11.86786178624856
$ws.Range("F9").Value = 38.2766472741239
$ws.Range("G9").Value = 48.51035149687426
$ws.Range("H9").Value = 17.96054065678829
$ws.Range("I9").Value = 23.85085537357724
$ws.Range("J9").Value = 9.205972277547968
$ws.Range("M9").Value = 32.27389189080321
$ws.Range("B10").Value = 11.11087104498522
$ws.Range("D10").Value = 7.992399254050572
$ws.Range("E10").Value = 11.8384570016868
$ws.Range("F10").Value = 38.69162054632616
$ws.Range("G10").Value = 49.92004175892642
$ws.Range("H10").Value = 18.04852271078686
$ws.Range("I10").Value = 24.80582862285623
$ws.Range("J10").Value = 9.12653119815162
$ws.Range("M10").Value = 34.14083581628081
$ws.Range("B11").Value = 11.40069902186998
$ws.Range("D11").Value = 8.02230108982735
$ws.Range("E11").Value = 11.82702849641727
$ws.Range("F11").Value = 38.8984355465104
$ws.Range("G11").Value = 50.57581542848442
$ws.Range("H11").Value = 18.09577415239711
$ws.Range("I11").Value = 25.23419270438996
$ws.Range("J11").Value = 9.092091459741027
$ws.Range("M11").Value = 34.95573371323803
$ws.Range("B12").Value = 11.50841786426322
$ws.Range("D12").Value = 8.033587358904832
$ws.Range("E12").Value = 11.82297986671853
$ws.Range("F12").Value = 38.97930953899657
$ws.Range("G12").Value = 50.82595957760634
$ws.Range("H12").Value = 18.11470623746566
$ws.Range("I12").Value = 25.39535805210319
$ws.Range("J12").Value = 9.079293140830437
$ws.Range("M12").Value = 35.25918984104987
$ws.Range("B13").Value = 11.48530976476248
$ws.Range("D13").Value = 8.031158343043151
$ws.Range("E13").Value = 11.82383941510549
$ws.Range("F13").Value = 38.96177878583146
$ws.Range("G13").Value = 50.77201000240821
$ws.Range("H13").Value = 18.11058269133033
$ws.Range("I13").Value = 25.36069742270268
$ws.Range("J13").Value = 9.082038683751742
$ws.Range("M13").Value = 35.1940658441995
$ws.Range("B14").Value = 11.40960209946328
$ws.Range("D14").Value = 8.023230372991586
$ws.Range("E14").Value = 11.82668982561478
$ws.Range("F14").Value = 38.90503807520692
$ws.Range("G14").Value = 50.59636005454512
$ws.Range("H14").Value = 18.09731089257746
$ws.Range("I14").Value = 25.24747381292164
$ws.Range("J14").Value = 9.091033664457097
$ws.Range("M14").Value = 34.98080296113932
$ws.Range("B15").Value = 11.36296296900496
$ws.Range("D15").Value = 8.018369376525964
$ws.Range("E15").Value = 11.82847210043459
$ws.Range("F15").Value = 38.87061462050481
$ws.Range("G15").Value = 50.48899807316653
$ws.Range("H15").Value = 18.08931679549069
$ws.Range("I15").Value = 25.17797977049789
$ws.Range("J15").Value = 9.096575001526524
$ws.Range("M15").Value = 34.84950047782856
$ws.Range("B16").Value = 11.09164806961416
$ws.Range("D16").Value = 7.990440079975575
$ws.Range("E16").Value = 11.83924299275825
$ws.Range("F16").Value = 38.6784647678716
$ws.Range("G16").Value = 49.87745422451905
$ws.Range("H16").Value = 18.04558017434333
$ws.Range("I16").Value = 24.77769679230438
$ws.Range("J16").Value = 9.128816009223222
$ws.Range("M16").Value = 34.08687210519246
$ws.Range("B17").Value = 10.9216331384878
$ws.Range("D17").Value = 7.973243488998534
$ws.Range("E17").Value = 11.84634878354496
$ws.Range("F17").Value = 38.56518411305367
$ws.Range("G17").Value = 49.50581442271176
$ws.Range("H17").Value = 18.02060076896718
$ws.Range("I17").Value = 24.53045525353447
$ws.Range("J17").Value = 9.149029131202401
$ws.Range("M17").Value = 33.61007827046247
$ws.Range("B18").Value = 10.82254815809527
$ws.Range("D18").Value = 7.96333021886257
$ws.Range("E18").Value = 11.85061925301404
$ws.Range("F18").Value = 38.50172890691206
$ws.Range("G18").Value = 49.29343778476781
$ws.Range("H18").Value = 18.00691375237964
$ws.Range("I18").Value = 24.38768790322158
$ws.Range("J18").Value = 9.160815102889277
$ws.Range("M18").Value = 33.33261178905103
$ws.Range("B19").Value = 10.78877854557276
$ws.Range("D19").Value = 7.959970016627967
$ws.Range("E19").Value = 11.85209669070552
$ws.Range("F19").Value = 38.48053718381007
$ws.Range("G19").Value = 49.22177602711771
$ws.Range("H19").Value = 18.00239641616339
$ws.Range("I19").Value = 24.33925862988065
$ws.Range("J19").Value = 9.164833126316925
$ws.Range("M19").Value = 33.23811816172233
$ws.Range("B20").Value = 10.93986611525645
$ws.Range("D20").Value = 7.975076417186814
$ws.Range("E20").Value = 11.84557338525213
$ws.Range("F20").Value = 38.57706724589312
$ws.Range("G20").Value = 49.54523515680938
$ws.Range("H20").Value = 18.0231894392344
$ws.Range("I20").Value = 24.55683386850562
$ws.Range("J20").Value = 9.146860865978608
$ws.Range("M20").Value = 33.66116921524561
$ws.Range("B21").Value = 11.43189478396071
$ws.Range("D21").Value = 8.025560032972948
$ws.Range("E21").Value = 11.82584502456535
$ws.Range("F21").Value = 38.92163511151067
$ws.Range("G21").Value = 50.6479055202831
$ws.Range("H21").Value = 18.10118095627297
$ws.Range("I21").Value = 25.28076007558967
$ws.Range("J21").Value = 9.088385026089053
$ws.Range("M21").Value = 35.04358394289395
$ws.Range("B22").Value = 11.74159913602485
$ws.Range("D22").Value = 8.058337627807965
$ws.Range("E22").Value = 11.8145776940517
$ws.Range("F22").Value = 39.16171607484322
$ws.Range("G22").Value = 51.3790319283717
$ws.Range("H22").Value = 18.15820741831463
$ws.Range("I22").Value = 25.74771400564447
$ws.Range("J22").Value = 9.051585172510082
$ws.Range("M22").Value = 35.91711498730143
$ws.Range("B23").Value = 11.57740345071119
$ws.Range("D23").Value = 8.040864291029274
$ws.Range("E23").Value = 11.82044281818832
$ws.Range("F23").Value = 39.03223222297667
$ws.Range("G23").Value = 50.98794623386522
$ws.Range("H23").Value = 18.12721795725363
$ws.Range("I23").Value = 25.49911134734642
$ws.Range("J23").Value = 9.071096559994837
$ws.Range("M23").Value = 35.45368995365267
$ws.Range("B24").Value = 10.93162716006264
$ws.Range("D24").Value = 7.974247833379971
$ws.Range("E24").Value = 11.84592336571461
$ws.Range("F24").Value = 38.57168967295031
$ws.Range("G24").Value = 49.52740902309515
$ws.Range("H24").Value = 18.02201700211462
$ws.Range("I24").Value = 24.54491002961143
$ws.Range("J24").Value = 9.147840624118343
$ws.Range("M24").Value = 33.63808142361886
$ws.Range("B25").Value = 10.18574457617505
$ws.Range("D25").Value = 7.901551197580489
$ws.Range("E25").Value = 11.88042390462244
$ws.Range("F25").Value = 38.14002509494034
$ws.Range("G25").Value = 48.00584574097135
$ws.Range("H25").Value = 17.93453995002665
$ws.Range("I25").Value = 23.49504671995851
$ws.Range("J25").Value = 9.236739835391012
$ws.Range("M25").Value = 31.55726515844424
